# Add a new "Proposed new matching mechanism" section at the end of the
# document: a page break, a Heading2 title, a quoted line, and a trailing
# blank paragraph.

$d = $word.ActiveDocument

# Position a collapsed range at the very end of the document's main story
# (right after the last paragraph mark), so the new content is appended
# after everything that's already there.
$endPos = $d.Content.End
$rng = $d.Range($endPos, $endPos)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$openQuote  = [char]0x201C
$enDash     = [char]0x2013
$rightQuote = [char]0x2019
$closeQuote = [char]0x201D

$quoteText = $openQuote + "You request my Classroom " + $enDash + " If you" + $rightQuote + "re bigger take it" + $closeQuote + " "

$xml  = "<w:p $wNs><w:r><w:br w:type=`"page`"/></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Proposed new matching mechanism</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/></w:pPr><w:r><w:t xml:space=`"preserve`">$quoteText</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/></w:pPr></w:p>"

$rng.InsertXML($xml)

Write-Output "Inserted proposed new matching mechanism section."
